# Jira defect1234: add a "country" column (with value "Australia") in front
# of the existing Email/password/repassword table on Sheet1, and normalize
# the "Email" header to lowercase "email".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column A; this shifts the existing
# Email/password/repassword columns one place to the right (B/C/D) and
# carries their column widths along with them.
$ws.Columns.Item(1).Insert()

# Column inserts don't repoint the worksheet's stored <hyperlink ref="..."/>
# (it keeps pointing at the old, now-empty A2), so rebuild it explicitly
# against the e-mail cell's new home in B2.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:gayatrik469@gmail.com")
$ws.Range("B2").Style = "Hyperlink"

# New leading "country" column.
$ws.Range("A1").Value = "country"
$ws.Range("A2").Value = "Australia"

# The old header was "Email"; the new layout uses lowercase "email".
$ws.Range("B1").Value = "email"

# Restore the authored selection (A3, single cell rather than A3:XFD24).
$ws.Range("A3").Select()
